$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows
$ws.Range("F3").Value = -6
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = 4
$ws.Range("F10").Value = 5
$ws.Range("F12").Value = 0
$ws.Range("F15").Value = -1
$ws.Range("F16").Value = -3
